$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: rename last column header from "DesDataHoracricao" to "DataHora"
$ws.Range("K1").Value = "DataHora"

# Row 2: fix the date/time text
$ws.Range("K2").Value = "18/03/2025 08:14"

# Row 3: description duplicate-fix ("Tese." typo corrected to match "Teste.") and date/time text
$ws.Range("J3").Value = "Teste.`nTeste."
$ws.Range("K3").Value = "18/03/2025 08:16"

# Row 4: fix the date/time text
$ws.Range("K4").Value = "18/03/2025 08:17"
